$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 9.462749333333333
$ws.Range("H2").Value = 28.388248
$ws.Range("I2").Value = 0.07254428564686972
$ws.Range("J2").Value = 0.07439525120506714
$ws.Range("M2").Value = 12.431794
$ws.Range("N2").Value = 37.295382
$ws.Range("O2").Value = 0.6267040910788743
$ws.Range("P2").Value = 0.7075740515758999
$ws.Range("Q2").Value = 117.6389503856373
$ws.Range("R2").Value = 1058.750553470736
$ws.Range("S2").Value = 0.04546380059928772
$ws.Range("T2").Value = 0.0526401493131762

$ws.Range("G3").Value = 9.462749333333333
$ws.Range("H3").Value = 28.388248
$ws.Range("I3").Value = 0.07254428564686972
$ws.Range("J3").Value = 0.07439525120506714
$ws.Range("O3").Value = 0.0264162940991436
$ws.Range("P3").Value = 0.0298250554119953
$ws.Range("Q3").Value = 4.95861628340711
$ws.Range("R3").Value = 44.62754655066399
$ws.Range("S3").Value = 0.001916351184859992
$ws.Range("T3").Value = 0.002218842489580438

$ws.Range("G4").Value = 9.462749333333333
$ws.Range("H4").Value = 28.388248
$ws.Range("I4").Value = 0.07254428564686972
$ws.Range("J4").Value = 0.07439525120506714
$ws.Range("M4").Value = 0.03915333333333333
$ws.Range("N4").Value = 0.11746
$ws.Range("O4").Value = 0.001973774193762771
$ws.Range("P4").Value = 0.002228470219130754
$ws.Range("Q4").Value = 0.3704981788977778
$ws.Range("R4").Value = 3.334483610079999
$ws.Range("S4").Value = 0.0001431860389147464
$ws.Range("T4").Value = 0.0001657876017552435

$ws.Range("G5").Value = 9.462749333333333
$ws.Range("H5").Value = 28.388248
$ws.Range("I5").Value = 0.07254428564686972
$ws.Range("J5").Value = 0.07439525120506714
$ws.Range("M5").Value = 6.8015495
$ws.Range("N5").Value = 13.603099
$ws.Range("O5").Value = 0.3428756056708687
$ws.Range("P5").Value = 0.2580802061075034
$ws.Range("Q5").Value = 64.36135799675867
$ws.Range("R5").Value = 386.168147980552
$ws.Range("S5").Value = 0.02487366587913096
$ws.Range("T5").Value = 0.01919994176442322

$ws.Range("G6").Value = 9.462749333333333
$ws.Range("H6").Value = 28.388248
$ws.Range("I6").Value = 0.07254428564686972
$ws.Range("J6").Value = 0.07439525120506714
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.04027333333333333
$ws.Range("N6").Value = 0.12082
$ws.Range("O6").Value = 0.002030234957350741
$ws.Range("P6").Value = 0.002292216685470609
$ws.Range("Q6").Value = 0.3810964581511111
$ws.Range("R6").Value = 3.429868123359999
$ws.Range("S6").Value = 0.0001472819446763125
$ws.Range("T6").Value = 0.0001705300361320324

$ws.Range("I7").Value = 0.3231336970688258
$ws.Range("J7").Value = 0.3313784449305509
$ws.Range("M7").Value = 12.431794
$ws.Range("N7").Value = 37.295382
$ws.Range("O7").Value = 0.6267040910788743
$ws.Range("P7").Value = 0.7075740515758999
$ws.Range("Q7").Value = 523.9986667240335
$ws.Range("R7").Value = 4715.9880005163
$ws.Range("S7").Value = 0.2025092099184748
$ws.Range("T7").Value = 0.2344747888844311

$ws.Range("I8").Value = 0.3231336970688258
$ws.Range("J8").Value = 0.3313784449305509
$ws.Range("O8").Value = 0.0264162940991436
$ws.Range("P8").Value = 0.0298250554119953
$ws.Range("S8").Value = 0.00853599477511368
$ws.Range("T8").Value = 0.009883380482394512

$ws.Range("I9").Value = 0.3231336970688258
$ws.Range("J9").Value = 0.3313784449305509
$ws.Range("M9").Value = 0.03915333333333333
$ws.Range("N9").Value = 0.11746
$ws.Range("O9").Value = 0.001973774193762771
$ws.Range("P9").Value = 0.002228470219130754
$ws.Range("Q9").Value = 1.650308432111111
$ws.Range("R9").Value = 14.852775889
$ws.Range("S9").Value = 0.0006377929524096052
$ws.Range("T9").Value = 0.0007384669957895932

$ws.Range("I10").Value = 0.3231336970688258
$ws.Range("J10").Value = 0.3313784449305509
$ws.Range("M10").Value = 6.8015495
$ws.Range("N10").Value = 13.603099
$ws.Range("O10").Value = 0.3428756056708687
$ws.Range("P10").Value = 0.2580802061075034
$ws.Range("Q10").Value = 286.6845179108917
$ws.Range("R10").Value = 1720.10710746535
$ws.Range("S10").Value = 0.1107946620951407
$ws.Range("T10").Value = 0.08552221736726053

$ws.Range("I11").Value = 0.3231336970688258
$ws.Range("J11").Value = 0.3313784449305509
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.04027333333333333
$ws.Range("N11").Value = 0.12082
$ws.Range("O11").Value = 0.002030234957350741
$ws.Range("P11").Value = 0.002292216685470609
$ws.Range("Q11").Value = 1.697516301444445
$ws.Range("R11").Value = 15.277646713
$ws.Range("S11").Value = 0.0006560373276871148
$ws.Range("T11").Value = 0.000759591200675112

$ws.Range("G12").Value = 32.300192
$ws.Range("H12").Value = 96.900576
$ws.Range("I12").Value = 0.2476229975407503
$ws.Range("J12").Value = 0.2539410918713864
$ws.Range("M12").Value = 12.431794
$ws.Range("N12").Value = 37.295382
$ws.Range("O12").Value = 0.6267040910788743
$ws.Range("P12").Value = 0.7075740515758999
$ws.Range("Q12").Value = 401.5493331044481
$ws.Range("R12").Value = 3613.943997940032
$ws.Range("S12").Value = 0.1551863456040022
$ws.Range("T12").Value = 0.1796821272370447

$ws.Range("G13").Value = 32.300192
$ws.Range("H13").Value = 96.900576
$ws.Range("I13").Value = 0.2476229975407503
$ws.Range("J13").Value = 0.2539410918713864
$ws.Range("O13").Value = 0.0264162940991436
$ws.Range("P13").Value = 0.0298250554119953
$ws.Range("Q13").Value = 16.92576357741866
$ws.Range("R13").Value = 152.331872196768
$ws.Range("S13").Value = 0.006541281928747972
$ws.Range("T13").Value = 0.00757380713644669

$ws.Range("G14").Value = 32.300192
$ws.Range("H14").Value = 96.900576
$ws.Range("I14").Value = 0.2476229975407503
$ws.Range("J14").Value = 0.2539410918713864
$ws.Range("M14").Value = 0.03915333333333333
$ws.Range("N14").Value = 0.11746
$ws.Range("O14").Value = 0.001973774193762771
$ws.Range("P14").Value = 0.002228470219130754
$ws.Range("Q14").Value = 1.264660184106667
$ws.Range("R14").Value = 11.38194165696
$ws.Range("S14").Value = 0.0004887518823281149
$ws.Range("T14").Value = 0.0005659001606489315

$ws.Range("G15").Value = 32.300192
$ws.Range("H15").Value = 96.900576
$ws.Range("I15").Value = 0.2476229975407503
$ws.Range("J15").Value = 0.2539410918713864
$ws.Range("M15").Value = 6.8015495
$ws.Range("N15").Value = 13.603099
$ws.Range("O15").Value = 0.3428756056708687
$ws.Range("P15").Value = 0.2580802061075034
$ws.Range("Q15").Value = 219.691354747504
$ws.Range("R15").Value = 1318.148128485024
$ws.Range("S15").Value = 0.08490388525982077
$ws.Range("T15").Value = 0.06553716932933186

$ws.Range("G16").Value = 32.300192
$ws.Range("H16").Value = 96.900576
$ws.Range("I16").Value = 0.2476229975407503
$ws.Range("J16").Value = 0.2539410918713864
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.04027333333333333
$ws.Range("N16").Value = 0.12082
$ws.Range("O16").Value = 0.002030234957350741
$ws.Range("P16").Value = 0.002292216685470609
$ws.Range("Q16").Value = 1.300836399146667
$ws.Range("R16").Value = 11.70752759232
$ws.Range("S16").Value = 0.0005027328658512078
$ws.Range("T16").Value = 0.0005820880079142168

$ws.Range("G17").Value = 9.736177999999999
$ws.Range("H17").Value = 19.472356
$ws.Range("I17").Value = 0.07464047213559308
$ws.Range("J17").Value = 0.0510299478916239
$ws.Range("M17").Value = 12.431794
$ws.Range("N17").Value = 37.295382
$ws.Range("O17").Value = 0.6267040910788743
$ws.Range("P17").Value = 0.7075740515758999
$ws.Range("Q17").Value = 121.038159243332
$ws.Range("R17").Value = 726.228955459992
$ws.Range("S17").Value = 0.0467774892474349
$ws.Range("T17").Value = 0.03610746698138338

$ws.Range("G18").Value = 9.736177999999999
$ws.Range("H18").Value = 19.472356
$ws.Range("I18").Value = 0.07464047213559308
$ws.Range("J18").Value = 0.0510299478916239
$ws.Range("O18").Value = 0.0264162940991436
$ws.Range("P18").Value = 0.0298250554119953
$ws.Range("Q18").Value = 5.101896823884665
$ws.Range("R18").Value = 30.61138094330799
$ws.Range("S18").Value = 0.00197172466363276
$ws.Range("T18").Value = 0.001521971023538916

$ws.Range("G19").Value = 9.736177999999999
$ws.Range("H19").Value = 19.472356
$ws.Range("I19").Value = 0.07464047213559308
$ws.Range("J19").Value = 0.0510299478916239
$ws.Range("M19").Value = 0.03915333333333333
$ws.Range("N19").Value = 0.11746
$ws.Range("O19").Value = 0.001973774193762771
$ws.Range("P19").Value = 0.002228470219130754
$ws.Range("Q19").Value = 0.3812038226266666
$ws.Range("R19").Value = 2.28722293576
$ws.Range("S19").Value = 0.0001473234377115028
$ws.Range("T19").Value = 0.0001137187191602781

$ws.Range("G20").Value = 9.736177999999999
$ws.Range("H20").Value = 19.472356
$ws.Range("I20").Value = 0.07464047213559308
$ws.Range("J20").Value = 0.0510299478916239
$ws.Range("M20").Value = 6.8015495
$ws.Range("N20").Value = 13.603099
$ws.Range("O20").Value = 0.3428756056708687
$ws.Range("P20").Value = 0.2580802061075034
$ws.Range("Q20").Value = 66.22109660781099
$ws.Range("R20").Value = 264.884386431244
$ws.Range("S20").Value = 0.02559239709105107
$ws.Range("T20").Value = 0.01316981946952545

$ws.Range("G21").Value = 9.736177999999999
$ws.Range("H21").Value = 19.472356
$ws.Range("I21").Value = 0.07464047213559308
$ws.Range("J21").Value = 0.0510299478916239
$ws.Range("K21").Value = 1
$ws.Range("L21").Value = 0.3333333333333333
$ws.Range("M21").Value = 0.04027333333333333
$ws.Range("N21").Value = 0.12082
$ws.Range("O21").Value = 0.002030234957350741
$ws.Range("P21").Value = 0.002292216685470609
$ws.Range("Q21").Value = 0.3921083419866666
$ws.Range("R21").Value = 2.35265005192
$ws.Range("S21").Value = 0.000151537695762845
$ws.Range("T21").Value = 0.000116971698015876

$ws.Range("G22").Value = 36.79199966666667
$ws.Range("H22").Value = 110.375999
$ws.Range("I22").Value = 0.2820585476079611
$ws.Range("J22").Value = 0.2892552641013719
$ws.Range("M22").Value = 12.431794
$ws.Range("N22").Value = 37.295382
$ws.Range("O22").Value = 0.6267040910788743
$ws.Range("P22").Value = 0.7075740515758999
$ws.Range("Q22").Value = 457.3905607040688
$ws.Range("R22").Value = 4116.515046336619
$ws.Range("S22").Value = 0.1767672457096747
$ws.Range("T22").Value = 0.2046695191598646

$ws.Range("G23").Value = 36.79199966666667
$ws.Range("H23").Value = 110.375999
$ws.Range("I23").Value = 0.2820585476079611
$ws.Range("J23").Value = 0.2892552641013719
$ws.Range("O23").Value = 0.0264162940991436
$ws.Range("P23").Value = 0.0298250554119953
$ws.Range("Q23").Value = 19.27953517732855
$ws.Range("R23").Value = 173.515816595957
$ws.Range("S23").Value = 0.007450941546789197
$ws.Range("T23").Value = 0.008627054280034751

$ws.Range("G24").Value = 36.79199966666667
$ws.Range("H24").Value = 110.375999
$ws.Range("I24").Value = 0.2820585476079611
$ws.Range("J24").Value = 0.2892552641013719
$ws.Range("M24").Value = 0.03915333333333333
$ws.Range("N24").Value = 0.11746
$ws.Range("O24").Value = 0.001973774193762771
$ws.Range("P24").Value = 0.002228470219130754
$ws.Range("Q24").Value = 1.440529426948889
$ws.Range("R24").Value = 12.96476484254
$ws.Range("S24").Value = 0.0005567198823988016
$ws.Range("T24").Value = 0.0006445967417767083

$ws.Range("G25").Value = 36.79199966666667
$ws.Range("H25").Value = 110.375999
$ws.Range("I25").Value = 0.2820585476079611
$ws.Range("J25").Value = 0.2892552641013719
$ws.Range("M25").Value = 6.8015495
$ws.Range("N25").Value = 13.603099
$ws.Range("O25").Value = 0.3428756056708687
$ws.Range("P25").Value = 0.2580802061075034
$ws.Range("Q25").Value = 250.2426069368169
$ws.Range("R25").Value = 1501.455641620901
$ws.Range("S25").Value = 0.0967109953457252
$ws.Range("T25").Value = 0.07465105817696237

$ws.Range("G26").Value = 36.79199966666667
$ws.Range("H26").Value = 110.375999
$ws.Range("I26").Value = 0.2820585476079611
$ws.Range("J26").Value = 0.2892552641013719
$ws.Range("K26").Value = 1
$ws.Range("L26").Value = 0.3333333333333333
$ws.Range("M26").Value = 0.04027333333333333
$ws.Range("N26").Value = 0.12082
$ws.Range("O26").Value = 0.002030234957350741
$ws.Range("P26").Value = 0.002292216685470609
$ws.Range("Q26").Value = 1.481736466575556
$ws.Range("R26").Value = 13.33562819918
$ws.Range("S26").Value = 0.0005726451233732608
$ws.Range("T26").Value = 0.0006630357427333723
